{"js": "// Remove \"modify and \" from the [R2] administrators requirement paragraph,\n// turning \"[R2] Administrators must be able to modify and cancel a taxi\n// driver's account.\" into \"[R2] Administrators must be able to cancel a\n// taxi driver's account.\" Word re-anchors its hidden \"_GoBack\" bookmark at\n// the most recent edit location, so it ends up split out between the two\n// halves of the sentence instead of trailing the whole paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"Administrators must be able to modify and cancel a taxi driver\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate target paragraph\");\n}\n\nconst results = target.search(\"modify and \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not locate 'modify and ' inside target paragraph\");\n}\n\n// Replacing the match with an empty string collapses the returned range to\n// the edit point, mirroring what Word does when text is deleted.\nconst matchRange = results.items[0];\nconst editPoint = matchRange.insertText(\"\", Word.InsertLocation.replace);\n\n// Move the hidden \"_GoBack\" bookmark to that collapsed edit point.\ncontext.document.deleteBookmark(\"_GoBack\");\neditPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that holds the [R2] \"modify and cancel\" requirement.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Administrators must be able to modify and cancel a taxi driver*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate target paragraph\"\n}\n\n$pr = $target.Range\n\n# Find \"modify and \" inside the paragraph; Find.Execute collapses the range\n# to the matched text so we can read off its start offset before deleting it.\n$find = $pr.Find\n$find.Text = \"modify and \"\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not locate 'modify and ' inside target paragraph\"\n}\n\n$cutStart = $pr.Start\n\n# Remove the matched text, collapsing the range to the edit point.\n$pr.Text = \"\"\n\n# Word re-anchors the hidden \"_GoBack\" bookmark at the most recent edit\n# location; reproduce that by moving/creating it at the collapsed point.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n$editPoint = $d.Range($cutStart, $cutStart)\n$d.Bookmarks.Add(\"_GoBack\", $editPoint)\n"}
